$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.044699999999994
$ws.Range("A12").Value = -21.57630000000001
$ws.Range("D14").Value = -7.463400000000002
$ws.Range("D26").Value = -8.433200000000005
$ws.Range("A27").Value = -21.5824
$ws.Range("D31").Value = -8.359900000000003
$ws.Range("A32").Value = -21.2123
$ws.Range("D35").Value = -8.339199999999998
$ws.Range("A36").Value = -20.4658
$ws.Range("D37").Value = -7.900299999999998
$ws.Range("A38").Value = -19.697
$ws.Range("D45").Value = -7.6662
$ws.Range("A46").Value = -21.82790000000001
$ws.Range("D52").Value = -7.6972
$ws.Range("A54").Value = -21.648
$ws.Range("A55").Value = -22.5819
$ws.Range("A56").Value = -22.2437
$ws.Range("D57").Value = -8.756299999999998
$ws.Range("A67").Value = -21.43859999999998
$ws.Range("A69").Value = -21.61039999999998
$ws.Range("A72").Value = -21.79079999999999
$ws.Range("D81").Value = -7.070399999999994
$ws.Range("A83").Value = -21.8551
$ws.Range("D83").Value = -8.382199999999997
$ws.Range("A86").Value = -22.0165
$ws.Range("A91").Value = -21.45410000000001
$ws.Range("A93").Value = -21.13219999999999
$ws.Range("A99").Value = -20.33749999999998
$ws.Range("D100").Value = -8.051400000000003
$ws.Range("D102").Value = -7.991900000000003
